# "user guide structure fix"
# Insert a new "Windage X" worksheet between "Windage " and "Horizontal surf",
# carrying the windage area/lever-arm reference table that used to live only
# as a paper guide. Existing sheets are left untouched apart from the tab
# selection moving to the freshly inserted sheet (Excel only ever keeps one
# sheet "active" / tabSelected, so adding+activating a new sheet naturally
# clears it on the previously active one).

$wb = $excel.ActiveWorkbook

# Worksheets.Add() with no args inserts right before the currently active
# sheet - "Horizontal surf" is active in the source workbook, so the new
# sheet lands exactly between "Windage " and "Horizontal surf", matching
# sheetId="3" / r:id="rId2" taking Horizontal surf's old slot.
$ws = $wb.Worksheets.Add()
$ws.Name = "Windage X"

# Column A + header, entered in the same order the source authored the
# shared-string table (A2:A5, then the "Area [m2]" header, then A6:A10) so
# the regenerated shared-strings pool lines up with the source workbook.
$ws.Range("A2").Value = "Выхлопная труба"
$ws.Range("A3").Value = "Рубка"
$ws.Range("A4").Value = "Кран кормовой "
$ws.Range("A5").Value = "Кран носовой"

# Header row (row 1): B1:D1 - reuses the same "X1 [m]" / "X2 [m]" headers
# already used on the Horizontal surf sheet, plus a new "Area [m2]" header.
# Centred, matching the header style used on the Horizontal surf sheet.
$ws.Range("B1").Value = "Area [m2]"
$ws.Range("C1").Value = "X1 [m]"
$ws.Range("D1").Value = "X2 [m]"
$ws.Range("B1:D1").HorizontalAlignment = -4108

$ws.Range("A6").Value = "Крышки и комингс кормового трюма"
$ws.Range("A7").Value = "Крышки и комингс носового трюма"
$ws.Range("A8").Value = "Ют"
$ws.Range("A9").Value = "Бак"
$ws.Range("A10").Value = "Надводный борт (выше 2,001)"

# Column widths: col A fits the longest label, B:D match the narrower
# numeric columns used elsewhere in the workbook.
$ws.Columns.Item(1).ColumnWidth = 35.5
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Range("C1:D1").EntireColumn.ColumnWidth = 9.17

# Match the source sheet's selection anchor.
$ws.Range("K26").Select() | Out-Null
